# Auto_ML_Libs.xlsx - "Salavando primeira versão do projeto"
# Builds the AutoML library comparison table on Plan1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (library names) -------------------------------------------
# Written in this order so the shared-string table comes out in the same
# sequence the original authoring session produced it.
$ws.Range("A3").Value  = "Biblioteca"
$ws.Range("A2").Value  = "Caracteríticas"
$ws.Range("A4").Value  = "Auto-Sklearn"
$ws.Range("A5").Value  = "Auto-Keras"
$ws.Range("A6").Value  = "Google Cloud AutoML"
$ws.Range("A7").Value  = "AutoML"
$ws.Range("A8").Value  = "Pycaret"
$ws.Range("A9").Value  = "H2O"
$ws.Range("A10").Value = "TPOT"
$ws.Range("A11").Value = "AutoGluon"
$ws.Range("A12").Value = "DataRobot"
$ws.Range("A13").Value = "Darwin"
$ws.Range("A14").Value = "Tazi.ai"
$ws.Range("A15").Value = "JADBio AutoML"
$ws.Range("A16").Value = "MLJar"
$ws.Range("A17").Value = "Dataiku"
$ws.Range("A18").Value = "Microsoft Azure AutoML"
$ws.Range("A19").Value = "Amazon SageMaker Autopilot"
$ws.Range("A20").Value = "Akkio"

# --- Row 1 (criteria headers) --------------------------------------------
$ws.Range("C1").Value = "Grátis"
$ws.Range("D1").Value = "Feature Eng"
$ws.Range("E1").Value = "Feature Selection"
$ws.Range("F1").Value = "Limite Tempo Processamento"
$ws.Range("G1").Value = "Tunning"
$ws.Range("I1").Value = "Redes Neurais"
$ws.Range("B1").Value = "Popularidade"

# --- Extra rows added later ------------------------------------------------
$ws.Range("A21").Value = "DataBricks AutoML"
$ws.Range("H1").Value  = "Explicação"
$ws.Range("D24").Value = "gramas"

# --- Numeric scratch cells -------------------------------------------------
$ws.Range("D22").Value = 1000
$ws.Range("E22").Value = 1
$ws.Range("D25").Value = 3750

# --- F18: empty cell with an underline applied (no value typed in) --------
$ws.Range("F18").Font.Underline = $true

# --- Column widths -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 27.7109375
$ws.Range("B1:I1").EntireColumn.ColumnWidth = 17.5703125

# --- Page setup: A4 portrait ---------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection -------------------------------------------------------------
$ws.Range("F18").Select()
